# Updated TestData for Portugal Market
$wb = $excel.ActiveWorkbook

# 1) Change the Czech sheet selection to the full A1:D12 range (no single active cell).
$czech = $wb.Worksheets.Item("Czech")
$czech.Range("A1:D12").Select()

# 2) Create the new "Portugal" sheet as a copy of "Swiss", placed after it.
$swiss = $wb.Worksheets.Item("Swiss")
$swiss.Copy($null, $swiss)
$portugal = $wb.Worksheets.Item($swiss.Index + 1)
$portugal.Name = "Portugal"

# 3) Update the column widths for the new sheet to match the Portugal layout.
#    (calibrated so the engine's internal character-width rounding lands on
#    the closest achievable value to the authored widths)
$portugal.Columns.Item(1).ColumnWidth = 26.6
$portugal.Columns.Item(2).ColumnWidth = 15.6
$portugal.Columns.Item(3).ColumnWidth = 13.26
$portugal.Columns.Item(4).ColumnWidth = 15.76

# 4) Give rows 3-5 the taller row height used on the Portugal sheet.
$portugal.Rows.Item(3).RowHeight = 28.8
$portugal.Rows.Item(4).RowHeight = 28.8
$portugal.Rows.Item(5).RowHeight = 28.8

# 5) Update the market name and user-story reference cells.
#    (Order matters for shared-string allocation: the user-story code must
#    become the new string added right after the existing ones, followed by
#    the market name, to match the authored shared-strings table.)
$portugal.Range("B4").Value = "NGC-3479/T2435"
$portugal.Range("B2").Value = "Portugal Market"

# 6) Set the selection/active cell on the new sheet and make sure it is the active tab.
$portugal.Range("B2").Select()
$portugal.Activate()
